$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1) - update 想去人数 (column F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1887
$ws1.Range("F6").Value  = 3187
$ws1.Range("F12").Value = 560
$ws1.Range("F13").Value = 409
$ws1.Range("F18").Value = 1654
$ws1.Range("F30").Value = 46
$ws1.Range("F35").Value = 86
$ws1.Range("F36").Value = 1655
$ws1.Range("F38").Value = 1912

# Sheet "全部类型" (sheetId 4) - update 想去人数 (column F) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1887
$ws4.Range("F6").Value  = 3187
$ws4.Range("F12").Value = 560
$ws4.Range("F14").Value = 409
$ws4.Range("F19").Value = 1654
$ws4.Range("F31").Value = 46
$ws4.Range("F38").Value = 86
$ws4.Range("F39").Value = 1655
$ws4.Range("F41").Value = 1912
